# Add new "send mail" workflow configuration rows to the Constants sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert 5 new rows before row 17 (shifts existing rows 17.. down to 22..)
$ws.Rows.Item(17).Resize(5).Insert()

# Populate the newly inserted rows with the mail server configuration
# (values are entered in this order to mirror the original authoring order)
$ws.Range("A17").Value = "ServerMail"
$ws.Range("B17").Value = "chkmailrelay.corp.sovos.local"

$ws.Range("A19").Value = "PortMail"

$ws.Range("A18").Value = "DefaultMailAdress"
$ws.Range("B18").Value = "uipathdevelopment@sovos.com"

$ws.Range("A20").Value = "DefaultMailName"
$ws.Range("B20").Value = "UiPath Development"

$ws.Range("B19").Value = 25

# Row 21 stays blank as a separator, but restore the row height that Insert()
# would otherwise leave at the sheet default for the whole new block.
$ws.Rows.Item(17).Resize(5).RowHeight = 14.25

# Update the active selection / view as seen after the edit
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("A19").Select()
